{"js": "// Apply the four text edits described in the diff using Word's Office.js API.\n\nasync function replaceOnce(context, searchText, newText) {\n  const results = context.document.body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\n// 1. Intro paragraph: add \"scientists, researchers, and\" before \"students\".\nawait replaceOnce(\n  context,\n  \"The goal of this course is for students to learn:\",\n  \"The goal of this course is for scientists, researchers, and students to learn:\"\n);\n\n// 2. \"By the end of the course...\" paragraph rewrite.\nawait replaceOnce(\n  context,\n  \"By the end of the course, the hope is that students will have mastered strategies that allow them to use the above-listed, freely-available and open-source tools for conducting their research in a reproducible fashion.\",\n  \"By the end of the course, the hope is that we will all have mastered strategies allowing us to use the above-listed, freely-available and open-source tools for conducting research in a reproducible fashion.\"\n);\n\n// 3. Fix typo \"Addtionally\" -> \"Additionally\".\nawait replaceOnce(context, \"Addtionally, packaging\", \"Additionally, packaging\");\n\n// 4. Rewrite the closing sentence about MS-Word / LaTeX / markdown-formatted articles.\nawait replaceOnce(\n  context,\n  \"might be expedient if you don't know LaTeX, and until the journals get their act together and start accepting markdown-formatted articles.\",\n  \"might be expedient if you don't know LaTeX. (Though I have heard rumors that some journals may start accepting accepting markdown-formatted articles, eventually).\"\n);\n", "ps1": "# Apply the four text edits described in the diff using the Word COM object model.\n#\n# We locate each target phrase with Range.Find.Execute (which collapses the\n# range to the matched text) and then assign Range.Text directly. Using\n# Find.Execute's own ReplaceWith parameter would trigger Word's smart-quotes\n# AutoCorrect and turn straight apostrophes into curly ones, so we avoid that\n# by writing the replacement text straight onto the found Range instead.\n\nfunction Replace-FirstMatch {\n    param(\n        $Doc,\n        [string]$OldText,\n        [string]$NewText\n    )\n\n    $rng = $Doc.Content\n    $found = $rng.Find.Execute($OldText)\n    if (-not $found) {\n        throw \"Could not find text: $OldText\"\n    }\n    $rng.Text = $NewText\n}\n\n$d = $word.ActiveDocument\n\n# 1. Intro paragraph: add \"scientists, researchers, and\" before \"students\".\nReplace-FirstMatch $d \"The goal of this course is for students to learn:\" \"The goal of this course is for scientists, researchers, and students to learn:\"\n\n# 2. \"By the end of the course...\" paragraph rewrite.\nReplace-FirstMatch $d \"By the end of the course, the hope is that students will have mastered strategies that allow them to use the above-listed, freely-available and open-source tools for conducting their research in a reproducible fashion.\" \"By the end of the course, the hope is that we will all have mastered strategies allowing us to use the above-listed, freely-available and open-source tools for conducting research in a reproducible fashion.\"\n\n# 3. Fix typo \"Addtionally\" -> \"Additionally\".\nReplace-FirstMatch $d \"Addtionally, packaging\" \"Additionally, packaging\"\n\n# 4. Rewrite the closing sentence about MS-Word / LaTeX / markdown-formatted articles.\nReplace-FirstMatch $d \"might be expedient if you don't know LaTeX, and until the journals get their act together and start accepting markdown-formatted articles.\" \"might be expedient if you don't know LaTeX. (Though I have heard rumors that some journals may start accepting accepting markdown-formatted articles, eventually).\"\n"}
